# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Summary of the change (derived from the target diff):
#  - "Valor Mora" total (E11) increases from 994994 to 1187266
#  - "Cant. Periodos" (F13) increases from 6 to 7 (one more period is added
#    for worker CANDIDA R. BOHORQUEZ ARDILA: period 2508)
#  - The detail table (rows 16-21) is re-sorted/updated:
#      * Worker JORGE HUMBERTO MORALES MENDOZA (15700482) keeps periods
#        2109 and 2110, now moved to the top (rows 16-17) with updated
#        "Salario Basico" values (G) of 4492340 and F value 131820 for
#        both periods.
#      * Worker CANDIDA R. BOHORQUEZ ARDILA (32684786) now has 5 periods
#        instead of 4: 2504, 2505, 2506, 2507, 2508 (rows 18-22), with a
#        brand-new row (22) added for period 2508.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for the new detail row by inserting a blank row at 22
#    (this pushes the two trailing signature rows down from 26/27 to
#    27/28, and the mergeCells shift automatically with it).
# ---------------------------------------------------------------------
$ws.Rows("22:22").Insert()

# Row 22 becomes the new "last row of the table" (it needs the special
# bottom-border formatting that row 21 currently carries), so copy that
# formatting down into row 22 first...
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)

# ...and then row 21 goes back to being a normal "middle" row, so give it
# the plain formatting that row 20 uses.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Update the summary fields above the table.
# ---------------------------------------------------------------------
$ws.Cells.Item(11, 5).Value = 1187266   # E11 - Valor Mora
$ws.Cells.Item(13, 6).Value = 7         # F13 - Cant. Periodos

# ---------------------------------------------------------------------
# 3) Rewrite the detail rows (16-22) with the refreshed data.
# ---------------------------------------------------------------------
# columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$rows = @(
    @{ Row = 16; Doc = "15700482"; Nombre = "JORGE HUMBERTO MORALES MENDOZA"; Periodo = "2109"; F = 131820; G = 4492340 },
    @{ Row = 17; Doc = "15700482"; Nombre = "JORGE HUMBERTO MORALES MENDOZA"; Periodo = "2110"; F = 131820; G = 4492340 },
    @{ Row = 18; Doc = "32684786"; Nombre = "CANDIDA R. BOHORQUEZ ARDILA";    Periodo = "2504"; F = 179694; G = 4806804 },
    @{ Row = 19; Doc = "32684786"; Nombre = "CANDIDA R. BOHORQUEZ ARDILA";    Periodo = "2505"; F = 179694; G = 4806804 },
    @{ Row = 20; Doc = "32684786"; Nombre = "CANDIDA R. BOHORQUEZ ARDILA";    Periodo = "2506"; F = 179694; G = 4806804 },
    @{ Row = 21; Doc = "32684786"; Nombre = "CANDIDA R. BOHORQUEZ ARDILA";    Periodo = "2507"; F = 192272; G = 4806804 },
    @{ Row = 22; Doc = "32684786"; Nombre = "CANDIDA R. BOHORQUEZ ARDILA";    Periodo = "2508"; F = 192272; G = 4806804 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = "CC"
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc
    $ws.Cells.Item($r.Row, 4).Value = $r.Nombre
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}

Write-Host "Workbook updated"
